$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "CSS Id" -> "VLJ #"
$ws.Range("C2").Value = "VLJ #"

# Judge Bernard Jones' CSS Id "BVAJONESB" -> "123" (rows 3-7)
$ws.Range("C3:C7").Value = "123"

# Judge Lauren Roth -> Stuart Huels, DSUSER -> 860 (rows 8-9)
$ws.Range("B8:B9").Value = "Huels, Stuart"
$ws.Range("C8:C9").Value = "860"

# Add a new, blank, formatted row at the bottom of the table (row 10),
# matching the look of the last existing data row.
$xlPasteFormats = -4122
$ws.Range("A9:I9").Copy()
$ws.Range("A10:I10").PasteSpecial($xlPasteFormats)
$ws.Rows(10).RowHeight = 17
$excel.CutCopyMode = $false
